$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update FromLocation (D2) value from "nyc" to "las"
$ws.Range("D2").Value = "las"

# Update the selected/active cell in the sheet view
$ws.Range("D13").Select()

# Update the workbook window width
$excel.Width = 18330
